# Apply a stock-quantity correction pass to the Companywise Stock Report.
#
# For a set of line items, the "Closing Qty" (column F) is revised, which in
# turn changes the line's "Closing Value" (column G = Unit Rate * Qty,
# column D holds the unit rate), each company's "Sub Total:" row (column B,
# sum of that company's line G values) and finally the last "Sub Total:" row
# plus the "Grand Total:" row (both equal to the sum of every company's
# subtotal).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> corrected quantity (column F)
$qtyUpdates = @{
    16  = 44
    51  = 69
    95  = 202
    109 = 3
    123 = 31
    124 = 0
    139 = 166
    145 = 54
    152 = 24
    155 = 48
    170 = 32
    174 = 0
    180 = 3
    192 = 13
    226 = 11
    228 = 40
    229 = 5
    251 = 4
    292 = 96
    312 = 187
    319 = 4
    324 = 118
    334 = 1
    339 = 81
    340 = 117
    345 = 56
    346 = 163
    351 = 15
    353 = 185
    361 = 3
    364 = 20
    381 = 187
    387 = 16
    393 = 18
    401 = 765
    414 = 8
    417 = 619
    423 = 29
    429 = 408
    431 = 201
    434 = 105
    437 = 110
    442 = 36
    443 = 27
    447 = 17
    467 = 50
    474 = 6
    490 = 125
    498 = 16
    516 = 174
    518 = 34
    520 = 58
    521 = 293
    522 = 278
    556 = 43
    563 = 15
    567 = 5
    568 = 38
    570 = 68
    571 = 33
    599 = 8
    701 = 3
    703 = 13
    733 = 8
    799 = 87
    802 = 41
    805 = 3
    843 = 230
    846 = 54
    847 = 173
    848 = 325
    852 = 78
    853 = 170
    855 = 379
    862 = 179
    863 = 335
    865 = 226
    877 = 25
    889 = 1
    895 = 4
    912 = 1546
    914 = 92
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Update quantity (F) and recompute line value (G = D * F) for each changed row.
# NOTE: this interop layer's parameterized "Value" getter can't be read back
# into a PowerShell variable directly, so reads use "Value2" (assignment via
# "Value"/"Value2" both work fine).
foreach ($row in $qtyUpdates.Keys) {
    $qty = $qtyUpdates[$row]
    $rate = $ws.Cells.Item($row, 4).Value2   # column D = unit rate
    $ws.Cells.Item($row, 6).Value2 = $qty    # column F = quantity
    $ws.Cells.Item($row, 7).Value2 = [math]::Round($rate * $qty, 10)  # column G = value
}

# Recompute every "Sub Total:" row in column A as the sum of column G for the
# block of item rows belonging to it (from just after the previous
# "Sub Total:"/section row down to itself).
$subtotalRows = New-Object System.Collections.Generic.List[int]
for ($r = 1; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 1).Value2
    if ($label -eq "Sub Total:") {
        $subtotalRows.Add($r) | Out-Null
    }
}

$prev = 0
$finalSubtotalRow = $subtotalRows[$subtotalRows.Count - 1]
foreach ($sr in $subtotalRows) {
    if ($sr -eq $finalSubtotalRow) {
        # The very last "Sub Total:" row is itself the sum of all the other
        # company subtotals (handled after this loop), not a plain item block.
        $prev = $sr
        continue
    }
    $sum = 0
    for ($r = $prev + 1; $r -lt $sr; $r++) {
        $gVal = $ws.Cells.Item($r, 7).Value2
        if ($gVal -is [double] -or $gVal -is [int]) {
            $sum += $gVal
        }
    }
    $ws.Cells.Item($sr, 2).Value2 = [math]::Round($sum, 10)
    $prev = $sr
}

# The last "Sub Total:" row and the "Grand Total:" row both equal the sum of
# every company subtotal computed above.
$grandSum = 0
for ($i = 0; $i -lt $subtotalRows.Count - 1; $i++) {
    $bVal = $ws.Cells.Item($subtotalRows[$i], 2).Value2
    if ($bVal -is [double] -or $bVal -is [int]) {
        $grandSum += $bVal
    }
}
$grandSum = [math]::Round($grandSum, 10)

$ws.Cells.Item($finalSubtotalRow, 2).Value2 = $grandSum

for ($r = $finalSubtotalRow + 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Value2 -eq "Grand Total:") {
        $ws.Cells.Item($r, 2).Value2 = $grandSum
        break
    }
}
